$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adiós")
$ws.Range("A5").Value = 42824.37045138889
$ws.Range("A5").NumberFormat = "mmm d yyyy hh:mm AM/PM"
